# Auto-generated edit script applying value updates described in the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 13750
$ws.Range("J21").Value = 13750
$ws.Range("L21").Value = 13750
$ws.Range("N21").Value = -14686
$ws.Range("H23").Value = 13750
$ws.Range("J23").Value = 13750
$ws.Range("L23").Value = 13750
$ws.Range("N23").Value = -14218
$ws.Range("H38").Value = 11995.091
$ws.Range("J38").Value = 5250
$ws.Range("L38").Value = 15750
$ws.Range("N38").Value = -16494
$ws.Range("H43").Value = 2626.25
$ws.Range("I43").Value = 1902.2
$ws.Range("K43").Value = 1902.2
$ws.Range("M43").Value = -1833.2
$ws.Range("H55").Value = 608
$ws.Range("I55").Value = 164
$ws.Range("J55").Value = 774.5
$ws.Range("K55").Value = 164
$ws.Range("L55").Value = 774.5
$ws.Range("M55").Value = 50
$ws.Range("N55").Value = -1202.5
$ws.Range("H58").Value = 6789.067
$ws.Range("H80").Value = 1204.1724
$ws.Range("J80").Value = 2453.8333
$ws.Range("L80").Value = 7361.499899999999
$ws.Range("N80").Value = -9357.499899999999
$ws.Range("H83").Value = 1204.1724
$ws.Range("J83").Value = 2453.8333
$ws.Range("L83").Value = 22084.4997
$ws.Range("N83").Value = -32068.4997
$ws.Range("H87").Value = 444467780
$ws.Range("J87").Value = 444467780
$ws.Range("L87").Value = 444467780
$ws.Range("N87").Value = -444470276
$ws.Range("H90").Value = 444467780
$ws.Range("J90").Value = 444467780
$ws.Range("L90").Value = 1333403340
$ws.Range("N90").Value = -1333415820
$ws.Range("H112").Value = 1985.4
$ws.Range("I112").Value = 1599.6666
$ws.Range("J112").Value = 2081.8333
$ws.Range("K112").Value = 4798.9998
$ws.Range("L112").Value = 6245.499899999999
$ws.Range("M112").Value = -3690.9998
$ws.Range("N112").Value = -8461.499899999999
$ws.Range("H116").Value = 15716.529
$ws.Range("I116").Value = 14242.556
$ws.Range("K116").Value = 14242.556
$ws.Range("M116").Value = -10800.556
$ws.Range("H137").Value = 4019.4
$ws.Range("I137").Value = 3939.8
$ws.Range("K137").Value = 11819.4
$ws.Range("M137").Value = -9269.400000000001
$ws.Range("H138").Value = 2764506.5
$ws.Range("I138").Value = 25379.8
$ws.Range("J138").Value = 3253636.2
$ws.Range("K138").Value = 76139.39999999999
$ws.Range("L138").Value = 9760908.600000001
$ws.Range("M138").Value = -70999.39999999999
$ws.Range("N138").Value = -9771188.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 110961.02
$ws.Range("I32").Value = 110961.02
$ws.Range("K32").Value = 110961.02
$ws.Range("M32").Value = -110674.02
$ws.Range("H122").Value = 7410333.5
$ws.Range("I122").Value = 12348712
$ws.Range("J122").Value = 2765.5
$ws.Range("K122").Value = 37046136
$ws.Range("L122").Value = 8296.5
$ws.Range("M122").Value = -37043686
$ws.Range("N122").Value = -13196.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I86").Value = 18468.572
$ws.Range("J86").Value = 4499.6665
$ws.Range("K86").Value = 18468.572
$ws.Range("L86").Value = 4499.6665
$ws.Range("M86").Value = -17345.572
$ws.Range("N86").Value = -6745.6665
$ws.Range("I89").Value = 18468.572
$ws.Range("J89").Value = 4499.6665
$ws.Range("K89").Value = 92342.86
$ws.Range("L89").Value = 22498.3325
$ws.Range("M89").Value = -86726.86
$ws.Range("N89").Value = -33730.3325
$ws.Range("H134").Value = 2671.3215
$ws.Range("I134").Value = 2304.0417
$ws.Range("K134").Value = 6912.125100000001
$ws.Range("M134").Value = -4377.125100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2416.1924
$ws.Range("I58").Value = 2266.158
$ws.Range("J58").Value = 2823.4285
$ws.Range("K58").Value = 2266.158
$ws.Range("L58").Value = 2823.4285
$ws.Range("M58").Value = -2063.158
$ws.Range("N58").Value = -3229.4285
$ws.Range("H132").Value = 16273.223
$ws.Range("I132").Value = 16810.654
$ws.Range("K132").Value = 50431.962
$ws.Range("M132").Value = -47901.962
$ws.Range("H136").Value = 2416.1924
$ws.Range("I136").Value = 2266.158
$ws.Range("J136").Value = 2823.4285
$ws.Range("K136").Value = 6798.474
$ws.Range("L136").Value = 8470.2855
$ws.Range("M136").Value = -4248.474
$ws.Range("N136").Value = -13570.2855
$ws.Range("H141").Value = 718436.25
$ws.Range("J141").Value = 718436.25
$ws.Range("L141").Value = 718436.25
$ws.Range("N141").Value = -728796.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 139593
$ws.Range("I7").Value = 248.5
$ws.Range("K7").Value = 745.5
$ws.Range("M7").Value = -633.5
$ws.Range("H92").Value = 390.62857
$ws.Range("J92").Value = 426.73334
$ws.Range("L92").Value = 1280.20002
$ws.Range("N92").Value = -3776.20002
$ws.Range("H113").Value = 29950
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 2500
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 22250
$ws.Range("J92").Value = 22250
$ws.Range("L92").Value = 22250
$ws.Range("N92").Value = -25994
$ws.Range("H122").Value = 7171.9546
$ws.Range("I122").Value = 5573.3076
$ws.Range("K122").Value = 16719.9228
$ws.Range("M122").Value = -14269.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3709.5
$ws.Range("I22").Value = 3889.8
$ws.Range("J22").Value = 3529.2
$ws.Range("K22").Value = 3889.8
$ws.Range("L22").Value = 3529.2
$ws.Range("M22").Value = -3594.8
$ws.Range("N22").Value = -4119.2
$ws.Range("H27").Value = 3709.5
$ws.Range("I27").Value = 3889.8
$ws.Range("J27").Value = 3529.2
$ws.Range("K27").Value = 3889.8
$ws.Range("L27").Value = 3529.2
$ws.Range("M27").Value = -3782.8
$ws.Range("N27").Value = -3743.2
$ws.Range("H46").Value = 6456.2144
$ws.Range("I46").Value = 4477
$ws.Range("J46").Value = 7555.778
$ws.Range("K46").Value = 4477
$ws.Range("L46").Value = 7555.778
$ws.Range("M46").Value = -4289
$ws.Range("N46").Value = -7931.778
$ws.Range("H61").Value = 5359.533
$ws.Range("I61").Value = 2056.2856
$ws.Range("J61").Value = 8249.875
$ws.Range("K61").Value = 2056.2856
$ws.Range("L61").Value = 8249.875
$ws.Range("M61").Value = -1854.2856
$ws.Range("N61").Value = -8653.875
$ws.Range("H113").Value = 5359.533
$ws.Range("I113").Value = 2056.2856
$ws.Range("J113").Value = 8249.875
$ws.Range("K113").Value = 2056.2856
$ws.Range("L113").Value = 8249.875
$ws.Range("M113").Value = 113.7143999999998
$ws.Range("N113").Value = -12589.875
$ws.Range("H121").Value = 136334
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("H132").Value = 1434195.2
$ws.Range("J132").Value = 2865413.5
$ws.Range("L132").Value = 8596240.5
$ws.Range("N132").Value = -8601300.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4316.6665
$ws.Range("I81").Value = 3180
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 6360
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = -5299
$ws.Range("N81").Value = -22122
$ws.Range("H84").Value = 4316.6665
$ws.Range("I84").Value = 3180
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 31800
$ws.Range("L84").Value = 100000
$ws.Range("M84").Value = -26496
$ws.Range("N84").Value = -110608
